# A new weekly price-report row for "Camote" (sweet potato) needs to be
# inserted into the consolidated data table. In the source system this
# lands in date order, which falls between the existing row 107 and the
# former row 108 - i.e. at worksheet row 108 - pushing every row below it
# down by one (107 unchanged, old 108 -> 109, old 201 -> 202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108; Excel shifts rows 108..201 down to 109..202
# and the sheet's used range (dimension) grows from R201 to R202 automatically.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = 'Vega Modelo de Temuco'
$ws.Range("C108").Value = 'La Araucanía'
$ws.Range("D108").Value = 45068
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = 100114002
$ws.Range("G108").Value = 'Camote'
$ws.Range("H108").Value = 'Sin especificar'
$ws.Range("I108").Value = 'Primera'
$ws.Range("J108").Value = 80
$ws.Range("K108").Value = 26000
$ws.Range("L108").Value = 26000
$ws.Range("M108").Value = 26000
$ws.Range("N108").Value = '$/caja 18 kilos'
$ws.Range("O108").Value = 'Perú'
$ws.Range("P108").Value = 1444
$ws.Range("Q108").Value = 18
$ws.Range("R108").Value = 'Hortaliza'
